$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
# Row 33
$ws.Range("H33").Value = 9870312
$ws.Range("J33").Value = 27780462
$ws.Range("L33").Value = 27780462
$ws.Range("N33").Value = -27780920

# Row 40
$ws.Range("H40").Value = 1982.6666
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 2024
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 2024
$ws.Range("M40").Value = -1725
$ws.Range("N40").Value = -2374

# Row 43
$ws.Range("H43").Value = 61906336
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4138

# Row 53
$ws.Range("H53").Value = 125000930
$ws.Range("I53").Value = 784
$ws.Range("K53").Value = 784
$ws.Range("M53").Value = -147

# Row 118
$ws.Range("H118").Value = 809.9
$ws.Range("I118").Value = 668.625
$ws.Range("K118").Value = 2005.875
$ws.Range("M118").Value = -348.875

# Row 132
$ws.Range("H132").Value = 1314.75
$ws.Range("I132").Value = 1243.3636
$ws.Range("K132").Value = 3730.0908
$ws.Range("M132").Value = -1200.0908

# Row 137
$ws.Range("H137").Value = 2785759.2
$ws.Range("I137").Value = 6054.25
$ws.Range("J137").Value = 6260390.5
$ws.Range("K137").Value = 18162.75
$ws.Range("L137").Value = 18781171.5
$ws.Range("M137").Value = -15612.75
$ws.Range("N137").Value = -18786271.5

# Row 138
$ws.Range("H138").Value = 5681.9
$ws.Range("I138").Value = 6121.593
$ws.Range("J138").Value = 4768.6924
$ws.Range("K138").Value = 18364.779
$ws.Range("L138").Value = 14306.0772
$ws.Range("M138").Value = -13224.779
$ws.Range("N138").Value = -24586.0772

$ws = $wb.Worksheets.Item(2)  # ARM
# Row 32
$ws.Range("H32").Value = 221117.44
$ws.Range("I32").Value = 271983.56
$ws.Range("K32").Value = 271983.56
$ws.Range("M32").Value = -271696.56

# Row 34
$ws.Range("H34").Value = 289666.66
$ws.Range("J34").Value = 500000
$ws.Range("L34").Value = 500000
$ws.Range("N34").Value = -500542

# Row 61
$ws.Range("H61").Value = 1696081.4
$ws.Range("I61").Value = 59679.895
$ws.Range("K61").Value = 59679.895
$ws.Range("M61").Value = -59467.895

# Row 74
$ws.Range("H74").Value = 933149
$ws.Range("I74").Value = 3971.353
$ws.Range("K74").Value = 3971.353
$ws.Range("M74").Value = -3097.353

# Row 77
$ws.Range("H77").Value = 933149
$ws.Range("I77").Value = 3971.353
$ws.Range("K77").Value = 19856.765
$ws.Range("M77").Value = -15488.765

# Row 136
$ws.Range("H136").Value = 1696081.4
$ws.Range("I136").Value = 59679.895
$ws.Range("K136").Value = 179039.685
$ws.Range("M136").Value = -176489.685

$ws = $wb.Worksheets.Item(3)  # BSM
# Row 99
$ws.Range("H99").Value = 9216
$ws.Range("I99").Value = 9581.647000000001
$ws.Range("K99").Value = 9581.647000000001
$ws.Range("M99").Value = -8083.647000000001

# Row 105
$ws.Range("H105").Value = 8348.843999999999
$ws.Range("I105").Value = 7035.391
$ws.Range("K105").Value = 7035.391
$ws.Range("M105").Value = -5288.391

$ws = $wb.Worksheets.Item(4)  # CRP
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# Row 31
$ws.Range("H31").Value = 2664.75
$ws.Range("I31").Value = 1985.5
$ws.Range("J31").Value = 5245.9
$ws.Range("K31").Value = 1985.5
$ws.Range("L31").Value = 5245.9
$ws.Range("M31").Value = -1690.5
$ws.Range("N31").Value = -5835.9

# Row 33
$ws.Range("H33").Value = 6487.1665
$ws.Range("J33").Value = 11937.5
$ws.Range("L33").Value = 11937.5
$ws.Range("N33").Value = -12695.5

# Row 34
$ws.Range("H34").Value = 2664.75
$ws.Range("I34").Value = 1985.5
$ws.Range("J34").Value = 5245.9
$ws.Range("K34").Value = 1985.5
$ws.Range("L34").Value = 5245.9
$ws.Range("M34").Value = -1783.5
$ws.Range("N34").Value = -5649.9

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 133
$ws.Range("H133").Value = 119369
$ws.Range("J133").Value = 119369
$ws.Range("L133").Value = 119369
$ws.Range("N133").Value = -124429

# Row 134
$ws.Range("H134").Value = 1709.6285
$ws.Range("I134").Value = 1411
$ws.Range("K134").Value = 4233
$ws.Range("M134").Value = -1698

$ws = $wb.Worksheets.Item(5)  # CUL
# Row 107
$ws.Range("H107").Value = 24391010
$ws.Range("J107").Value = 47619852
$ws.Range("L107").Value = 142859556
$ws.Range("N107").Value = -142863396

# Row 113
$ws.Range("H113").Value = 1520
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 6699.999899999999
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -11039.9999

# Row 115
$ws.Range("H115").Value = 3348
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 3348
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 10044
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -12394

# Row 117
$ws.Range("H117").Value = 166668180
$ws.Range("J117").Value = 166668180
$ws.Range("L117").Value = 500004540
$ws.Range("N117").Value = -500011424

# Row 119
$ws.Range("H119").Value = 10175.637
$ws.Range("J119").Value = 23996.666
$ws.Range("L119").Value = 71989.99800000001
$ws.Range("N119").Value = -81665.99800000001

# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# Row 121
$ws.Range("H121").Value = 4169.2
$ws.Range("I121").Value = 2077.75
$ws.Range("J121").Value = 4692.0625
$ws.Range("K121").Value = 6233.25
$ws.Range("L121").Value = 14076.1875
$ws.Range("M121").Value = -4923.25
$ws.Range("N121").Value = -16696.1875

# Row 122
$ws.Range("H122").Value = 5557333.5
$ws.Range("I122").Value = 8333889
$ws.Range("K122").Value = 75005001
$ws.Range("M122").Value = -75002551

$ws = $wb.Worksheets.Item(6)  # GSM
# Row 96
$ws.Range("H96").Value = 47420.332
$ws.Range("J96").Value = 47420.332
$ws.Range("L96").Value = 47420.332
$ws.Range("N96").Value = -52912.332

# Row 100
$ws.Range("H100").Value = 37199.6
$ws.Range("I100").Value = 29998
$ws.Range("J100").Value = 39000
$ws.Range("K100").Value = 29998
$ws.Range("L100").Value = 39000
$ws.Range("M100").Value = -28916
$ws.Range("N100").Value = -41164

# Row 102
$ws.Range("H102").Value = 33335128
$ws.Range("I102").Value = 41668410
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 41668410
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = -41666788
$ws.Range("N102").Value = -5243

# Row 122
$ws.Range("H122").Value = 2110.5
$ws.Range("I122").Value = 1861.6
$ws.Range("K122").Value = 5584.799999999999
$ws.Range("M122").Value = -3134.799999999999

# Row 132
$ws.Range("H132").Value = 1711198.4
$ws.Range("I132").Value = 1577.25
$ws.Range("J132").Value = 8549683
$ws.Range("K132").Value = 4731.75
$ws.Range("L132").Value = 25649049
$ws.Range("M132").Value = -2201.75
$ws.Range("N132").Value = -25654109

$ws = $wb.Worksheets.Item(7)  # LTW
# Row 7
$ws.Range("H7").Value = 9547.111000000001
$ws.Range("J7").Value = 12073.333
$ws.Range("L7").Value = 12073.333
$ws.Range("N7").Value = -12297.333

# Row 61
$ws.Range("H61").Value = 17511.25
$ws.Range("J61").Value = 16000
$ws.Range("L61").Value = 16000
$ws.Range("N61").Value = -16404

# Row 87
$ws.Range("H87").Value = 333383330
$ws.Range("J87").Value = 333383330
$ws.Range("L87").Value = 333383330
$ws.Range("N87").Value = -333385576

# Row 90
$ws.Range("H90").Value = 333383330
$ws.Range("J90").Value = 333383330
$ws.Range("L90").Value = 1000149990
$ws.Range("N90").Value = -1000161222

# Row 93
$ws.Range("H93").Value = 1412.8
$ws.Range("I93").Value = 1412.8
$ws.Range("K93").Value = 1412.8
$ws.Range("M93").Value = -164.8

# Row 113
$ws.Range("H113").Value = 17511.25
$ws.Range("J113").Value = 16000
$ws.Range("L113").Value = 16000
$ws.Range("N113").Value = -20340

# Row 122
$ws.Range("H122").Value = 3019.5625
$ws.Range("I122").Value = 2751.1667
$ws.Range("J122").Value = 3824.75
$ws.Range("K122").Value = 8253.500100000001
$ws.Range("L122").Value = 11474.25
$ws.Range("M122").Value = -5803.500100000001
$ws.Range("N122").Value = -16374.25

# Row 126
$ws.Range("H126").Value = 9547.111000000001
$ws.Range("J126").Value = 12073.333
$ws.Range("L126").Value = 36219.999
$ws.Range("N126").Value = -41159.999

$ws = $wb.Worksheets.Item(8)  # WVR
# Row 107
$ws.Range("H107").Value = 2383649.5
$ws.Range("I107").Value = 2482
$ws.Range("K107").Value = 7446
$ws.Range("M107").Value = -5526

# Row 122
$ws.Range("H122").Value = 1999.5834
$ws.Range("I122").Value = 1785
$ws.Range("K122").Value = 5355
$ws.Range("M122").Value = -2905

# Row 126
$ws.Range("H126").Value = 8336920
$ws.Range("I126").Value = 11367066
$ws.Range("K126").Value = 34101198
$ws.Range("M126").Value = -34098728

# Row 132
$ws.Range("H132").Value = 3483.2222
$ws.Range("I132").Value = 3071
$ws.Range("K132").Value = 9213
$ws.Range("M132").Value = -6683

# Row 136
$ws.Range("H136").Value = 638.6070999999999
$ws.Range("I136").Value = 547.375
$ws.Range("K136").Value = 1642.125
$ws.Range("M136").Value = 907.875
